# Update the task list: mark a few more tasks as assigned to Arthur, flip
# two of them to "Done", and flip one to "In Progress" (highlighted red).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("12. Feature to search the most "red" news"): already Done,
# just assign it to Arthur.
$ws.Range("D12").Value = "Arthur"

# Row 14 ("14. Fix bug with dissappearing tags"): mark Done + assign Arthur.
$ws.Range("C14").Value = "Done"
$ws.Range("C14").Font.Color = 5287936   # green FF00B050, same as other "Done" cells
$ws.Range("D14").Value = "Arthur"

# Row 15 ("15. Add button check all checkboxes"): mark Done + assign Arthur.
$ws.Range("C15").Value = "Done"
$ws.Range("C15").Font.Color = 5287936   # green FF00B050
$ws.Range("D15").Value = "Arthur"

# Row 17 ("17. Add links to sites"): mark In Progress (highlighted red) + assign Arthur.
$ws.Range("C17").Value = "In Progress"
$ws.Range("C17").Font.Color = 255       # red FFFF0000
$ws.Range("D17").Value = "Arthur"

# Leave the selection where the author left it while editing.
$ws.Range("F7").Select() | Out-Null
